$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update last_edited_time (shared by rows 6-13, which all shared the same
# original timestamp string) so the shared string table is updated in place
# rather than a new distinct string being appended.
$ws.Range("D6:D13").Value = "2024-08-12T02:00:00.000Z"

# Update numeric metrics for row 6 (Tháng 8)
$ws.Range("W6").Value = 100250000
$ws.Range("AA6").Value = 110000000
$ws.Range("AE6").Value = 210250000
$ws.Range("AH6").Value = 178250000
$ws.Range("AK6").Value = 26
$ws.Range("AN6").Value = 32000000
$ws.Range("AQ6").Value = 217250000
